# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" columns for the most recently
# handed-off file (89d1f092-2166-438a-82b5-a3c48562134e) in each
# language sheet, and rolls the newest of those timestamps up into the
# "Latest HO Xliff Generate Date" column on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: row 7 is the 89d1f092-....md entry; column H is "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-08-25 22:41:15"

# de-de: row 7 is the 89d1f092-....md entry; column H is "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-08-25 22:41:20"

# Overview: row 7 is the 89d1f092-....md entry; column G is "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-08-25 22:41:20"
